$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Slide 26: merge three runs ("Dans  la construction des modeles" /
#    ", cet " / "algorithme travaille de maniere sequentielle. ")
#    into a single run. The visible text does not change, only the
#    run (formatting-run) boundaries collapse into one <a:r>.
# -----------------------------------------------------------------
$s26 = $p.Slides.Item(26)
$sh26 = $s26.Shapes.Item(2)
$tf26 = $sh26.TextFrame
$tr26 = $tf26.TextRange
$full26 = $tr26.Text

$marker = [char]0x00E8  # 'e' with grave accent, used below to build the search text safely
$search = "Dans  la construction des mod" + $marker + "les, cet algorithme travaille de mani" + $marker + "re s" + [char]0x00E9 + "quentielle. "

$startIdx0 = $full26.IndexOf($search)
if ($startIdx0 -ge 0) {
    $startIdx1 = $startIdx0 + 1   # COM Characters() is 1-based
    $len = $search.Length
    $sub = $tr26.Characters($startIdx1, $len)
    $sub.Text = $sub.Text
}

# -----------------------------------------------------------------
# 2) Swap slides 28 and 29 (the "VUE DE l'ALGORITHME" slide and the
#    "Premier resultat" slide trade places in the deck order).
# -----------------------------------------------------------------
$p.Slides.Item(28).MoveTo(29)
